# Update column G (K - strikeouts) values for rows 2-32 with regenerated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 2
    3  = 3
    4  = 4
    5  = 3
    6  = 2
    7  = 3
    8  = 10
    9  = 5
    10 = 4
    11 = 6
    12 = 6
    13 = 7
    14 = 10
    15 = 6
    16 = 9
    17 = 4
    18 = 6
    19 = 8
    20 = 3
    21 = 5
    22 = 7
    23 = 5
    24 = 3
    25 = 3
    26 = 10
    27 = 8
    28 = 5
    29 = 5
    30 = 6
    31 = 5
    32 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
